# Femacal de La Calera - Espárragos: add a new weekly record.
# A new row of data is inserted at row 5 (pushing the existing rows 5-46
# down to 6-47), and the new row is populated with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; everything below shifts down by one.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44881
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 300000000
$ws.Range("G5").Value = "Espárragos"
$ws.Range("H5").Value = "Verde"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 680
$ws.Range("K5").Value = 1400
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = 1400
$ws.Range("N5").Value = "$/kilo"
$ws.Range("O5").Value = "Provincia de Quillota"
$ws.Range("P5").Value = 1400
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
